$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting the existing row 4
# (Singapore - Hougang vs Geylang) down to row 5.
$ws.Rows.Item(4).Insert()

# Fill the newly inserted row 4 with the Estonia Meistriliiga match.
$ws.Cells.Item(4, 1).Value = "vaVGVu14"
$ws.Cells.Item(4, 2).Value = "24/11/2024"
$ws.Cells.Item(4, 3).Value = "07:30"
$ws.Cells.Item(4, 4).Value = "ESTONIA - MEISTRILIIGA"
$ws.Cells.Item(4, 5).Value = "Viimsi JK"
$ws.Cells.Item(4, 6).Value = "Tallinna Kalev"
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(4, 8).Value = 4.1
$ws.Cells.Item(4, 9).Value = 1.5
$ws.Cells.Item(4, 10).Value = 5
$ws.Cells.Item(4, 11).Value = 2.32
$ws.Cells.Item(4, 12).Value = 2
$ws.Cells.Item(4, 15).Value = 1.17
$ws.Cells.Item(4, 16).Value = 3.78
$ws.Cells.Item(4, 17).Value = 1.6
$ws.Cells.Item(4, 18).Value = 2.07
$ws.Cells.Item(4, 19).Value = 1.31
$ws.Cells.Item(4, 20).Value = 3.26
$ws.Cells.Item(4, 21).Value = 1.78
$ws.Cells.Item(4, 22).Value = 1.99
$ws.Cells.Item(4, 23).Value = 13.5
$ws.Cells.Item(4, 24).Value = 25
$ws.Cells.Item(4, 25).Value = 13.5
$ws.Cells.Item(4, 26).Value = 70
$ws.Cells.Item(4, 27).Value = 37
$ws.Cells.Item(4, 28).Value = 35
$ws.Cells.Item(4, 29).Value = 13.5
$ws.Cells.Item(4, 30).Value = 7.2
$ws.Cells.Item(4, 31).Value = 13.5
$ws.Cells.Item(4, 32).Value = 50
$ws.Cells.Item(4, 33).Value = 300
$ws.Cells.Item(4, 34).Value = 6.8
$ws.Cells.Item(4, 35).Value = 6.6
$ws.Cells.Item(4, 36).Value = 7
$ws.Cells.Item(4, 37).Value = 9
$ws.Cells.Item(4, 38).Value = 9.5
$ws.Cells.Item(4, 39).Value = 18.5
$ws.Cells.Item(4, 40).Value = 6.8
$ws.Cells.Item(4, 41).Value = 28
$ws.Cells.Item(4, 42).Value = 32
$ws.Cells.Item(4, 43).Value = 175
$ws.Cells.Item(4, 44).Value = 175
$ws.Cells.Item(4, 45).Value = 400
$ws.Cells.Item(4, 46).Value = 3.1
$ws.Cells.Item(4, 47).Value = 7.6
$ws.Cells.Item(4, 48).Value = 65
$ws.Cells.Item(4, 49).Value = 3.4
$ws.Cells.Item(4, 50).Value = 7
$ws.Cells.Item(4, 51).Value = 15.5
$ws.Cells.Item(4, 52).Value = 20
$ws.Cells.Item(4, 53).Value = 45
$ws.Cells.Item(4, 54).Value = 200

# Add new row 6: Turkey Super Lig match.
$ws.Cells.Item(6, 1).Value = "xppFPDhg"
$ws.Cells.Item(6, 2).Value = "24/11/2024"
$ws.Cells.Item(6, 3).Value = "07:30"
$ws.Cells.Item(6, 4).Value = "TURKEY - SUPER LIG"
$ws.Cells.Item(6, 5).Value = "Sivasspor"
$ws.Cells.Item(6, 6).Value = "Kasimpasa"
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 3.4
$ws.Cells.Item(6, 9).Value = 2.25
$ws.Cells.Item(6, 10).Value = 3.5
$ws.Cells.Item(6, 11).Value = 2.2
$ws.Cells.Item(6, 12).Value = 2.88
$ws.Cells.Item(6, 13).Value = 1.04
$ws.Cells.Item(6, 14).Value = 13
$ws.Cells.Item(6, 15).Value = 1.22
$ws.Cells.Item(6, 16).Value = 4
$ws.Cells.Item(6, 17).Value = 1.75
$ws.Cells.Item(6, 18).Value = 2.05
$ws.Cells.Item(6, 19).Value = 1.36
$ws.Cells.Item(6, 20).Value = 3
$ws.Cells.Item(6, 21).Value = 1.62
$ws.Cells.Item(6, 22).Value = 2.2
$ws.Cells.Item(6, 23).Value = 11
$ws.Cells.Item(6, 24).Value = 17
$ws.Cells.Item(6, 25).Value = 11
$ws.Cells.Item(6, 26).Value = 29
$ws.Cells.Item(6, 27).Value = 23
$ws.Cells.Item(6, 28).Value = 29
$ws.Cells.Item(6, 29).Value = 12
$ws.Cells.Item(6, 30).Value = 6.5
$ws.Cells.Item(6, 31).Value = 13
$ws.Cells.Item(6, 32).Value = 41
$ws.Cells.Item(6, 33).Value = 151
$ws.Cells.Item(6, 34).Value = 9
$ws.Cells.Item(6, 35).Value = 12
$ws.Cells.Item(6, 36).Value = 9.5
$ws.Cells.Item(6, 37).Value = 21
$ws.Cells.Item(6, 38).Value = 17
$ws.Cells.Item(6, 39).Value = 23
$ws.Cells.Item(6, 40).Value = 5
$ws.Cells.Item(6, 41).Value = 17
$ws.Cells.Item(6, 42).Value = 23
$ws.Cells.Item(6, 43).Value = 51
$ws.Cells.Item(6, 44).Value = 67
$ws.Cells.Item(6, 45).Value = 151
$ws.Cells.Item(6, 46).Value = 3
$ws.Cells.Item(6, 47).Value = 7.5
$ws.Cells.Item(6, 48).Value = 51
$ws.Cells.Item(6, 49).Value = 4.5
$ws.Cells.Item(6, 50).Value = 12
$ws.Cells.Item(6, 51).Value = 21
$ws.Cells.Item(6, 52).Value = 41
$ws.Cells.Item(6, 53).Value = 51
$ws.Cells.Item(6, 54).Value = 126
$ws.Cells.Item(6, 55).Value = 251
$ws.Cells.Item(6, 56).Value = 301

# Add new row 7: Turkey 1. Lig match.
$ws.Cells.Item(7, 1).Value = "d2jjMXa3"
$ws.Cells.Item(7, 2).Value = "24/11/2024"
$ws.Cells.Item(7, 3).Value = "07:30"
$ws.Cells.Item(7, 4).Value = "TURKEY - 1. LIG"
$ws.Cells.Item(7, 5).Value = "Erzurumspor"
$ws.Cells.Item(7, 6).Value = "Karagumruk"
$ws.Cells.Item(7, 7).Value = 2.15
$ws.Cells.Item(7, 8).Value = 3.3
$ws.Cells.Item(7, 9).Value = 3.25
$ws.Cells.Item(7, 10).Value = 3
$ws.Cells.Item(7, 11).Value = 2.05
$ws.Cells.Item(7, 12).Value = 4
$ws.Cells.Item(7, 13).Value = 1.06
$ws.Cells.Item(7, 14).Value = 10
$ws.Cells.Item(7, 15).Value = 1.33
$ws.Cells.Item(7, 16).Value = 3.25
$ws.Cells.Item(7, 17).Value = 2.1
$ws.Cells.Item(7, 18).Value = 1.7
$ws.Cells.Item(7, 19).Value = 1.5
$ws.Cells.Item(7, 20).Value = 2.5
$ws.Cells.Item(7, 21).Value = 1.91
$ws.Cells.Item(7, 22).Value = 1.8
$ws.Cells.Item(7, 23).Value = 7
$ws.Cells.Item(7, 24).Value = 9.5
$ws.Cells.Item(7, 25).Value = 9.5
$ws.Cells.Item(7, 26).Value = 21
$ws.Cells.Item(7, 27).Value = 19
$ws.Cells.Item(7, 28).Value = 34
$ws.Cells.Item(7, 29).Value = 8.5
$ws.Cells.Item(7, 30).Value = 6.5
$ws.Cells.Item(7, 31).Value = 17
$ws.Cells.Item(7, 32).Value = 51
$ws.Cells.Item(7, 33).Value = 351
$ws.Cells.Item(7, 34).Value = 8.5
$ws.Cells.Item(7, 35).Value = 15
$ws.Cells.Item(7, 36).Value = 12
$ws.Cells.Item(7, 37).Value = 34
$ws.Cells.Item(7, 38).Value = 29
$ws.Cells.Item(7, 39).Value = 41
$ws.Cells.Item(7, 40).Value = 4
$ws.Cells.Item(7, 41).Value = 13
$ws.Cells.Item(7, 42).Value = 26
$ws.Cells.Item(7, 43).Value = 41
$ws.Cells.Item(7, 44).Value = 67
$ws.Cells.Item(7, 45).Value = 201
$ws.Cells.Item(7, 46).Value = 2.5
$ws.Cells.Item(7, 47).Value = 8.5
$ws.Cells.Item(7, 48).Value = 67
$ws.Cells.Item(7, 49).Value = 5
$ws.Cells.Item(7, 50).Value = 19
$ws.Cells.Item(7, 51).Value = 29
$ws.Cells.Item(7, 52).Value = 67
$ws.Cells.Item(7, 53).Value = 101
$ws.Cells.Item(7, 54).Value = 251
$ws.Cells.Item(7, 55).Value = 126
$ws.Cells.Item(7, 56).Value = 126

